$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 163.6875
$ws.Range("I33").Value = 163.6875
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 163.6875
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = $null
$ws.Range("N33").Value = 65.3125
$ws.Range("H88").Value = 17177.834
$ws.Range("I88").Value = 590
$ws.Range("J88").Value = 25471.75
$ws.Range("K88").Value = 590
$ws.Range("L88").Value = 25471.75
$ws.Range("M88").Value = -184
$ws.Range("N88").Value = -26283.75
$ws.Range("H91").Value = 17177.834
$ws.Range("I91").Value = 590
$ws.Range("J91").Value = 25471.75
$ws.Range("K91").Value = 590
$ws.Range("L91").Value = 25471.75
$ws.Range("M91").Value = 814
$ws.Range("N91").Value = -28279.75
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = $null
$ws.Range("N98").Value = 0
$ws.Range("H113").Value = 8624.25
$ws.Range("I113").Value = 7874.75
$ws.Range("J113").Value = 8999
$ws.Range("K113").Value = 7874.75
$ws.Range("L113").Value = 8999
$ws.Range("M113").Value = -4620.75
$ws.Range("N113").Value = -15507
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = $null
$ws.Range("N122").Value = 0
$ws.Range("H137").Value = 3183.5454
$ws.Range("I137").Value = 2814.9092
$ws.Range("J137").Value = 3552.182
$ws.Range("K137").Value = 8444.7276
$ws.Range("L137").Value = 10656.546
$ws.Range("M137").Value = -5894.7276
$ws.Range("N137").Value = -15756.546

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 27950
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 27950
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = $null
$ws.Range("M43").Value = 27950
$ws.Range("N43").Value = -28576
$ws.Range("H46").Value = 2826
$ws.Range("I46").Value = 2076
$ws.Range("J46").Value = 3576
$ws.Range("K46").Value = 2076
$ws.Range("L46").Value = 3576
$ws.Range("M46").Value = -1757
$ws.Range("N46").Value = -4214
$ws.Range("H74").Value = 4999
$ws.Range("I74").Value = 4999
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 4999
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -4125
$ws.Range("H77").Value = 4999
$ws.Range("I77").Value = 4999
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 24995
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -20627
$ws.Range("H121").Value = 33085
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 33085
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 33085
$ws.Range("N121").Value = -36579
$ws.Range("H132").Value = 2012
$ws.Range("I132").Value = 2012
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6036
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3506
$ws.Range("H139").Value = 49997.2
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 49997.2
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 49997.2
$ws.Range("N139").Value = -60277.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 46749.75
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 46749.75
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 46749.75
$ws.Range("N81").Value = -48871.75
$ws.Range("H84").Value = 46749.75
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 46749.75
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 140249.25
$ws.Range("N84").Value = -150857.25
$ws.Range("H86").Value = 6461.2
$ws.Range("I86").Value = 5835.3335
$ws.Range("J86").Value = 7400
$ws.Range("K86").Value = 5835.3335
$ws.Range("L86").Value = 7400
$ws.Range("M86").Value = -4712.3335
$ws.Range("N86").Value = -9646
$ws.Range("H89").Value = 6461.2
$ws.Range("I89").Value = 5835.3335
$ws.Range("J89").Value = 7400
$ws.Range("K89").Value = 29176.6675
$ws.Range("L89").Value = 37000
$ws.Range("M89").Value = -23560.6675
$ws.Range("N89").Value = -48232
$ws.Range("H134").Value = 13375
$ws.Range("I134").Value = 13375
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 40125
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -37590

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 3692.5557
$ws.Range("I12").Value = 1058.25
$ws.Range("J12").Value = 5800
$ws.Range("K12").Value = 1058.25
$ws.Range("L12").Value = 5800
$ws.Range("M12").Value = -888.25
$ws.Range("N12").Value = -6140
$ws.Range("H16").Value = 551.1111
$ws.Range("I16").Value = 475.33334
$ws.Range("J16").Value = 702.6667
$ws.Range("K16").Value = 475.33334
$ws.Range("L16").Value = 702.6667
$ws.Range("M16").Value = -188.33334
$ws.Range("N16").Value = -1276.6667
$ws.Range("H68").Value = 28571.428
$ws.Range("I68").Value = 20000
$ws.Range("H71").Value = 28571.428
$ws.Range("I71").Value = 20000
$ws.Range("H86").Value = 8040.875
$ws.Range("I86").Value = 6119.25
$ws.Range("J86").Value = 9962.5
$ws.Range("K86").Value = 6119.25
$ws.Range("L86").Value = 9962.5
$ws.Range("M86").Value = -4996.25
$ws.Range("N86").Value = -12208.5
$ws.Range("H89").Value = 8040.875
$ws.Range("I89").Value = 6119.25
$ws.Range("J89").Value = 9962.5
$ws.Range("K89").Value = 30596.25
$ws.Range("L89").Value = 49812.5
$ws.Range("M89").Value = -24980.25
$ws.Range("N89").Value = -61044.5
$ws.Range("H113").Value = 551.1111
$ws.Range("I113").Value = 475.33334
$ws.Range("J113").Value = 702.6667
$ws.Range("K113").Value = 475.33334
$ws.Range("L113").Value = 702.6667
$ws.Range("M113").Value = 1694.66666
$ws.Range("N113").Value = -5042.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 3515
$ws.Range("I130").Value = 1030
$ws.Range("J130").Value = 6000
$ws.Range("K130").Value = 3090
$ws.Range("L130").Value = 18000
$ws.Range("M130").Value = 1930
$ws.Range("N130").Value = -28040
$ws.Range("H131").Value = 1316.3334
$ws.Range("I131").Value = 971.5
$ws.Range("J131").Value = 1661.1666
$ws.Range("K131").Value = 2914.5
$ws.Range("L131").Value = 4983.4998
$ws.Range("M131").Value = 2125.5
$ws.Range("N131").Value = -15063.4998
$ws.Range("H136").Value = 2962.5
$ws.Range("I136").Value = 2962.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 8887.5
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -3787.5
$ws.Range("H137").Value = 1629
$ws.Range("I137").Value = 1629
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 4887
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = 213
$ws.Range("N137").Value = $null
$ws.Range("H138").Value = 2061.2856
$ws.Range("I138").Value = 1904.8334
$ws.Range("J138").Value = 3000
$ws.Range("K138").Value = 5714.5002
$ws.Range("L138").Value = 9000
$ws.Range("M138").Value = -574.5002000000004
$ws.Range("N138").Value = -19280
$ws.Range("H140").Value = 2484.375
$ws.Range("I140").Value = 1982.5714
$ws.Range("J140").Value = 5997
$ws.Range("K140").Value = 5947.7142
$ws.Range("L140").Value = 17991
$ws.Range("M140").Value = -767.7142000000003
$ws.Range("N140").Value = -28351

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 40000
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 40000
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 40000
$ws.Range("M63").Value = $null
$ws.Range("N63").Value = -41372
$ws.Range("H66").Value = 40000
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 40000
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 120000
$ws.Range("M66").Value = $null
$ws.Range("N66").Value = -126864
$ws.Range("H70").Value = 333333340
$ws.Range("I70").Value = 333333340
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 333333340
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -333333070
$ws.Range("H73").Value = 333333340
$ws.Range("I73").Value = 333333340
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 333333340
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -333332404
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = $null
$ws.Range("N80").Value = 0
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = $null
$ws.Range("N83").Value = 0
$ws.Range("H126").Value = 6552.8
$ws.Range("I126").Value = 3609.2
$ws.Range("J126").Value = 12440
$ws.Range("K126").Value = 10827.6
$ws.Range("L126").Value = 37320
$ws.Range("M126").Value = -8357.599999999999
$ws.Range("N126").Value = -42260
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = $null
$ws.Range("N127").Value = 0
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3387
$ws.Range("I61").Value = 3387
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3387
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3185
$ws.Range("H68").Value = 30125.5
$ws.Range("I68").Value = 3500.6667
$ws.Range("J68").Value = 110000
$ws.Range("K68").Value = 3500.6667
$ws.Range("L68").Value = 110000
$ws.Range("M68").Value = -2751.6667
$ws.Range("N68").Value = -111498
$ws.Range("H71").Value = 30125.5
$ws.Range("I71").Value = 3500.6667
$ws.Range("J71").Value = 110000
$ws.Range("K71").Value = 17503.3335
$ws.Range("L71").Value = 550000
$ws.Range("M71").Value = -13759.3335
$ws.Range("N71").Value = -557488
$ws.Range("H103").Value = 39450
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 39450
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 39450
$ws.Range("N103").Value = -41794
$ws.Range("H113").Value = 3387
$ws.Range("I113").Value = 3387
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3387
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1217
